# Aspen fires plot workbook update
# - Adds a second worksheet ("version 02") that is a trimmed copy of Sheet1:
#   the three intermediate "Quaking aspen suitability" rows (SSP245, SSP585,
#   Delta SSP245) are removed and the remaining delta-suitability row is
#   relabeled "2041-2070 Aspen suitability change (SSP585)".
# - Cleans up a stray cell format on Sheet1 (B18) that no longer needs its
#   own border-applying style.
# - Updates the active selections left on each sheet.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: normalize B18's style (was using a now-redundant bordered xf) ---
$sheet1.Range("B18").Borders.LineStyle = -4142

# --- Create "version 02" as a duplicate of Sheet1, placed right after it ---
$sheet1.Copy([System.Reflection.Missing]::Value, $sheet1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "version 02"

# Remove the SSP245 / SSP585 / Delta(SSP245) rows from the aspen-suitability
# block, leaving just "historic" and the (relabeled) delta-585 row.
$ws2.Rows("4:6").Delete()

# Relabel the remaining suitability-change row and match its original height.
$ws2.Range("C3").Value = "2041-2070 Aspen suitability change (SSP585) "
$ws2.Range("D3").Value = "delta585"
$ws2.Rows(4).RowHeight = 36

# Mirror the same B-cell style cleanup on the copied sheet (former B18, now
# row 15 after the row deletion above).
$ws2.Range("B15").Borders.LineStyle = -4142

# --- Selections ---
$ws2.Range("C4").Select()
$sheet1.Activate()
$sheet1.Range("C8").Select()
